# Minor changes: mark test cases 14-19 (rows 20-25) on Sheet1 as "Passed"
# instead of "Failed", and leave the active cell selection on H11.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()

# Rows 20-25, column J ("Status") go from "Failed" to "Passed"
foreach ($r in 20..25) {
    $ws1.Cells.Item($r, 10).Value = "Passed"
}

# Update the cell selection shown in the sheet view
$ws1.Range("H11").Select()

$wb.Save()
